$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.984.43"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.917.38"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.439"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.26%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "33.52"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("D15").Value = "3.397.66"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "61.017.95"
$ws.Range("E16").Value = "  +0.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.69"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "2.913.39"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "435.30"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.676"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.11"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.59"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.79"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.27"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.92%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.99"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -3.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.58"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "0.0₃0868"
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("E35").Value = "  +0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.62"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  -1.20%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.58"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.04"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.287"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.32%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "375.81"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D45").Value = "2.686.96"
$ws.Range("E45").Value = "  -0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.83"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.106"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.124"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.67%  "
